$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45957
$ws.Range("B2").Value = 87
$ws.Range("C2").Value = 77.42
$ws.Range("D2").Value = 76.93000000000001
$ws.Range("E2").Value = 77.12
$ws.Range("F2").Value = 79.59999999999999
$ws.Range("G2").Value = 103.5
$ws.Range("H2").Value = 110.49
$ws.Range("I2").Value = 152.68
$ws.Range("J2").Value = 124.09
$ws.Range("K2").Value = 85.19
$ws.Range("L2").Value = 51
$ws.Range("M2").Value = 26.32
$ws.Range("N2").Value = 29.77
$ws.Range("O2").Value = 21.77
$ws.Range("P2").Value = 20.03
$ws.Range("Q2").Value = 42.8
$ws.Range("R2").Value = 70.76000000000001
$ws.Range("S2").Value = 116.14
$ws.Range("T2").Value = 138.5
$ws.Range("U2").Value = 145.26
$ws.Range("V2").Value = 157.13
$ws.Range("W2").Value = 137.7
$ws.Range("X2").Value = 114.78
$ws.Range("Y2").Value = 107.43
$ws.Range("Z2").Value = 89.73
$ws.Range("AB2").Value = 129.26
$ws.Range("AD2").Value = 147.42
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 141.88
$ws.Range("AG2").Value = "0h-16h"
